$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.134.63"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").Value = "1.660.45"
$ws.Range("E3").Value = "  +4.06%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'215.60"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").Value = "'0.0615"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.893.67"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").Value = "1.667.41"
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "'64.99"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "'240.95"
$ws.Range("E17").Value = "  +6.20%  "
$ws.Range("D18").Value = "27.137.86"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "'7.87"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +5.47%  "
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("D25").Value = "'145.89"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'15.86"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "1.533.01"
$ws.Range("E32").Value = "  +6.13%  "
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").Value = "'3.06"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("E35").Value = "  +8.68%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.575"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "'0.894"
$ws.Range("E38").Value = "  +9.32%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'2.28"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").Value = "'66.29"
$ws.Range("E43").Value = "  +9.92%  "
$ws.Range("D44").Value = "1.800.06"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("D45").Value = "'0.772"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").Value = "'0.919"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'90.48"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "'0.0980"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").Value = "'0.0504"
$ws.Range("E51").Value = "  +0.63%  "
